$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells so numeric-looking strings (e.g. "585.80",
# "0.0000174") are kept exactly as text instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.684.30"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.430.90"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.80"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.04"
$ws.Range("E6").Value = "  +5.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +6.16%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.430.09"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.99"
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.024.32"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.46"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.587.12"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000174"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.398.74"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.88"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.30"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.535"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +6.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.85"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.41"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.10"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.28"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.870"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.67"
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.44"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.714.19"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.29"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0691"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.13"
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "339.27"
$ws.Range("E47").Value = "  +11.47%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.94"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.43"
$ws.Range("E50").Value = "  +7.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  +4.00%  "
